$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> { D = new price text (or $null if unchanged); E = new volume text (or $null if unchanged) }
$updates = @(
    @{ Row = 2; D = '43.660.62'; E = '  -0.09%  ' }
    @{ Row = 3; D = '2.291.57'; E = '  -0.93%  ' }
    @{ Row = 4; D = $null; E = '  +0.09%  ' }
    @{ Row = 5; D = '96.49'; E = '  +4.00%  ' }
    @{ Row = 6; D = '268.04'; E = '  -0.16%  ' }
    @{ Row = 7; D = $null; E = '  -1.62%  ' }
    @{ Row = 8; D = $null; E = '  +0.02%  ' }
    @{ Row = 9; D = '0.611'; E = '  -1.45%  ' }
    @{ Row = 10; D = '45.93'; E = '  +2.35%  ' }
    @{ Row = 11; D = $null; E = '  -0.04%  ' }
    @{ Row = 12; D = '7.83'; E = '  -2.57%  ' }
    @{ Row = 13; D = '0.106'; E = '  +0.47%  ' }
    @{ Row = 14; D = '2.634.29'; E = '  -0.68%  ' }
    @{ Row = 15; D = '15.16'; E = '  -0.72%  ' }
    @{ Row = 16; D = '0.849'; E = '  -0.77%  ' }
    @{ Row = 17; D = '2.297.20'; E = '  +0.06%  ' }
    @{ Row = 18; D = '43.615.49'; E = '  -0.43%  ' }
    @{ Row = 19; D = $null; E = '  +2.52%  ' }
    @{ Row = 20; D = '6.22'; E = '  -0.71%  ' }
    @{ Row = 21; D = '72.41'; E = '  +1.74%  ' }
    @{ Row = 22; D = '2.48'; E = '  +8.96%  ' }
    @{ Row = 23; D = '232.90'; E = '  -2.35%  ' }
    @{ Row = 24; D = '9.20'; E = '  -4.16%  ' }
    @{ Row = 25; D = $null; E = '  -0.11%  ' }
    @{ Row = 26; D = '2.54'; E = '  +1.57%  ' }
    @{ Row = 27; D = '11.20'; E = '  +0.39%  ' }
    @{ Row = 28; D = '3.48'; E = '  +2.42%  ' }
    @{ Row = 29; D = '40.40'; E = '  +3.69%  ' }
    @{ Row = 30; D = $null; E = '  -0.96%  ' }
    @{ Row = 31; D = '175.66'; E = '  +1.50%  ' }
    @{ Row = 32; D = '21.83'; E = '  -3.81%  ' }
    @{ Row = 33; D = '0.0893'; E = '  +0.59%  ' }
    @{ Row = 34; D = '5.37'; E = '  -2.34%  ' }
    @{ Row = 35; D = $null; E = '  -0.29%  ' }
    @{ Row = 36; D = '0.108'; E = '  -1.96%  ' }
    @{ Row = 37; D = '0.0354'; E = '  +1.46%  ' }
    @{ Row = 38; D = '4.35'; E = '  -2.76%  ' }
    @{ Row = 39; D = '3.39'; E = '  +1.09%  ' }
    @{ Row = 40; D = '0.240'; E = '  +2.13%  ' }
    @{ Row = 41; D = $null; E = '  -1.70%  ' }
    @{ Row = 42; D = '12.32'; E = '  +0.93%  ' }
    @{ Row = 43; D = $null; E = '  +0.77%  ' }
    @{ Row = 44; D = '64.65'; E = '  +5.83%  ' }
    @{ Row = 45; D = $null; E = '  -1.19%  ' }
    @{ Row = 46; D = $null; E = '  -4.22%  ' }
    @{ Row = 47; D = $null; E = '  +0.42%  ' }
    @{ Row = 48; D = '97.62'; E = '  -2.65%  ' }
    @{ Row = 49; D = $null; E = '  +0.29%  ' }
    @{ Row = 50; D = '2.514.36'; E = '  -0.40%  ' }
    @{ Row = 51; D = '0.429'; E = '  -0.54%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        # Prices like "96.49" or "9.20" parse as numbers (losing trailing
        # zeros / padding) unless the cell is forced to Text first. Multi-dot
        # values (e.g. "43.660.62") are never numeric so this is harmless for them.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

Write-Output "Updated $($updates.Count) rows"
